$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 38
$ws.Range("H38").Value = 100
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# row 39
$ws.Range("H39").Value = 47.142857
$ws.Range("J39").Value = 79
$ws.Range("L39").Value = 237
$ws.Range("N39").Value = -829

# row 40
$ws.Range("H40").Value = 2425
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

# row 61
$ws.Range("H61").Value = 970
$ws.Range("I61").Value = 970
$ws.Range("K61").Value = 2910
$ws.Range("M61").Value = -2738

# row 76
$ws.Range("H76").Value = 166675580
$ws.Range("I76").Value = 166675580
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 166675580
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -166675265
$ws.Range("N76").ClearContents()

# row 79
$ws.Range("H79").Value = 166675580
$ws.Range("I79").Value = 166675580
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 166675580
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -166674488
$ws.Range("N79").ClearContents()

# row 80
$ws.Range("H80").Value = 1160.25
$ws.Range("I80").Value = 444.75
$ws.Range("J80").Value = 1875.75
$ws.Range("K80").Value = 1334.25
$ws.Range("L80").Value = 5627.25
$ws.Range("M80").Value = -336.25
$ws.Range("N80").Value = -7623.25

# row 82
$ws.Range("H82").Value = 1236.8
$ws.Range("I82").Value = 1236.8
$ws.Range("K82").Value = 3710.4
$ws.Range("M82").Value = -3304.4

# row 83
$ws.Range("H83").Value = 1160.25
$ws.Range("I83").Value = 444.75
$ws.Range("J83").Value = 1875.75
$ws.Range("K83").Value = 4002.75
$ws.Range("L83").Value = 16881.75
$ws.Range("M83").Value = 989.25
$ws.Range("N83").Value = -26865.75

# row 85
$ws.Range("H85").Value = 1236.8
$ws.Range("I85").Value = 1236.8
$ws.Range("K85").Value = 3710.4
$ws.Range("M85").Value = -2306.4

# row 92
$ws.Range("H92").Value = 1916.625
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# row 98
$ws.Range("H98").Value = 1575.4138
$ws.Range("I98").Value = 1613.8214
$ws.Range("K98").Value = 1613.8214
$ws.Range("M98").Value = -115.8214

# row 99
$ws.Range("H99").Value = 333342850
$ws.Range("I99").Value = 999
$ws.Range("J99").Value = 500013760
$ws.Range("K99").Value = 2997
$ws.Range("L99").Value = 1500041280
$ws.Range("M99").Value = -1499
$ws.Range("N99").Value = -1500044276

# row 122
$ws.Range("H122").Value = 1575.4138
$ws.Range("I122").Value = 1613.8214
$ws.Range("K122").Value = 4841.4642
$ws.Range("M122").Value = -2391.4642

# row 138
$ws.Range("H138").Value = 3875.3171
$ws.Range("I138").Value = 1037.0667
$ws.Range("J138").Value = 5512.769
$ws.Range("K138").Value = 3111.2001
$ws.Range("L138").Value = 16538.307
$ws.Range("M138").Value = 2028.7999
$ws.Range("N138").Value = -26818.307


$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 1408514.5
$ws.Range("I32").Value = 630.54095
$ws.Range("K32").Value = 630.54095
$ws.Range("M32").Value = -343.54095

# row 45
$ws.Range("H45").Value = 2367.7273
$ws.Range("I45").Value = 1535.1428
$ws.Range("J45").Value = 3824.75
$ws.Range("K45").Value = 1535.1428
$ws.Range("L45").Value = 3824.75
$ws.Range("M45").Value = -1158.1428
$ws.Range("N45").Value = -4578.75

# row 63
$ws.Range("H63").Value = 4678.143
$ws.Range("J63").Value = 6344.3335
$ws.Range("L63").Value = 6344.3335
$ws.Range("N63").Value = -7716.3335

# row 66
$ws.Range("H66").Value = 4678.143
$ws.Range("J66").Value = 6344.3335
$ws.Range("L66").Value = 31721.6675
$ws.Range("N66").Value = -38585.6675

# row 102
$ws.Range("H102").Value = 5269.4614
$ws.Range("I102").Value = 5291.9165
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 5291.9165
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -3669.9165
$ws.Range("N102").Value = -8244

# row 132
$ws.Range("H132").Value = 729694.75
$ws.Range("I132").Value = 1214543
$ws.Range("J132").Value = 56294.445
$ws.Range("K132").Value = 3643629
$ws.Range("L132").Value = 168883.335
$ws.Range("M132").Value = -3641099
$ws.Range("N132").Value = -173943.335


$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 8097.826
$ws.Range("I94").Value = 2950.8667
$ws.Range("J94").Value = 17748.375
$ws.Range("K94").Value = 2950.8667
$ws.Range("L94").Value = 17748.375
$ws.Range("M94").Value = -2499.8667
$ws.Range("N94").Value = -18650.375

# row 134
$ws.Range("H134").Value = 1040120.1
$ws.Range("I134").Value = 1194706.6
$ws.Range("J134").Value = 9543.833000000001
$ws.Range("K134").Value = 3584119.8
$ws.Range("L134").Value = 28631.499
$ws.Range("M134").Value = -3581584.8
$ws.Range("N134").Value = -33701.499


$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 5029.857
$ws.Range("I31").Value = 1834
$ws.Range("K31").Value = 1834
$ws.Range("M31").Value = -1539

# row 34
$ws.Range("H34").Value = 5029.857
$ws.Range("I34").Value = 1834
$ws.Range("K34").Value = 1834
$ws.Range("M34").Value = -1632

# row 52
$ws.Range("H52").Value = 72250
$ws.Range("J52").Value = 72250
$ws.Range("L52").Value = 72250
$ws.Range("N52").Value = -72838

# row 132
$ws.Range("H132").Value = 5377.385
$ws.Range("I132").Value = 3724.147
$ws.Range("K132").Value = 11172.441
$ws.Range("M132").Value = -8642.440999999999

# row 137
$ws.Range("H137").Value = 69999.664
$ws.Range("J137").Value = 69999.664
$ws.Range("L137").Value = 69999.664
$ws.Range("N137").Value = -80199.664


$ws = $wb.Worksheets.Item("CUL")
# row 50
$ws.Range("H50").Value = 999.6667
$ws.Range("I50").Value = 999.6667
$ws.Range("K50").Value = 2999.0001
$ws.Range("M50").Value = -2518.0001

# row 53
$ws.Range("H53").Value = 999.6667
$ws.Range("I53").Value = 999.6667
$ws.Range("K53").Value = 2999.0001
$ws.Range("M53").Value = -2518.0001

# row 122
$ws.Range("H122").Value = 119566.12
$ws.Range("I122").Value = 260.18182
$ws.Range("K122").Value = 2341.63638
$ws.Range("M122").Value = 108.3636200000001


$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 6580.457
$ws.Range("I102").Value = 5602.7896
$ws.Range("K102").Value = 5602.7896
$ws.Range("M102").Value = -3980.7896

# row 107
$ws.Range("H107").Value = 1802.1428
$ws.Range("I107").Value = 524.2
$ws.Range("J107").Value = 4997
$ws.Range("K107").Value = 524.2
$ws.Range("L107").Value = 4997
$ws.Range("M107").Value = 1395.8
$ws.Range("N107").Value = -8837

# row 122
$ws.Range("H122").Value = 4355.1113
$ws.Range("I122").Value = 2420.8572
$ws.Range("J122").Value = 11125
$ws.Range("K122").Value = 7262.571599999999
$ws.Range("L122").Value = 33375
$ws.Range("M122").Value = -4812.571599999999
$ws.Range("N122").Value = -38275


$ws = $wb.Worksheets.Item("LTW")
# row 68
$ws.Range("H68").Value = 2614.95
$ws.Range("I68").Value = 2393.2666
$ws.Range("J68").Value = 3280
$ws.Range("K68").Value = 2393.2666
$ws.Range("L68").Value = 3280
$ws.Range("M68").Value = -1644.2666
$ws.Range("N68").Value = -4778

# row 71
$ws.Range("H71").Value = 2614.95
$ws.Range("I71").Value = 2393.2666
$ws.Range("J71").Value = 3280
$ws.Range("K71").Value = 11966.333
$ws.Range("L71").Value = 16400
$ws.Range("M71").Value = -8222.332999999999
$ws.Range("N71").Value = -23888

# row 122
$ws.Range("H122").Value = 5642.6
$ws.Range("I122").Value = 3112
$ws.Range("K122").Value = 9336
$ws.Range("M122").Value = -6886


$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 17771.143
$ws.Range("J62").Value = 18162.25
$ws.Range("L62").Value = 18162.25
$ws.Range("N62").Value = -19410.25

# row 65
$ws.Range("H65").Value = 17771.143
$ws.Range("J65").Value = 18162.25
$ws.Range("L65").Value = 90811.25
$ws.Range("N65").Value = -97051.25

# row 107
$ws.Range("H107").Value = 531.6429000000001
$ws.Range("I107").Value = 531.6429000000001
$ws.Range("K107").Value = 1594.9287
$ws.Range("M107").Value = 325.0712999999998

# row 122
$ws.Range("H122").Value = 2665.4443
$ws.Range("I122").Value = 2711.125
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 8133.375
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -5683.375
$ws.Range("N122").Value = -11800

# row 132
$ws.Range("H132").Value = 6693.193
$ws.Range("I132").Value = 4404
$ws.Range("K132").Value = 13212
$ws.Range("M132").Value = -10682

